$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Tukey HSD results (Diff/Lower/Upper/q-value/p-value/H0) after
# adding Chauvenet criterion filtering + std-based recompute.
$data = @(
    @{ Row=2; C=0.5547905909351704; D=0.5250135446218048; E=0.5845676372485359; F=81.74969316723494; G=0.001; H="False" }
    @{ Row=3; C=0.01831129196337701; D=-0.008445705850044921; E=0.04506828977679894; F=3.002757530800678; G=0.459844638217055; H="True" }
    @{ Row=4; C=0.01512096774193527; D=-0.01159416253387981; E=0.04183609801775034; F=2.483481626266732; G=0.6850479383718692; H="True" }
    @{ Row=5; C=0.002340354396522791; D=-0.02435181491503168; E=0.02903252370807726; F=0.3847126055834243; G=0.9; H="True" }
    @{ Row=6; C=0.01204819277108412; D=-0.01464855186585174; E=0.03874493740801997; F=1.980169017112511; G=0.9; H="True" }
    @{ Row=7; C=0.09827777777778073; D=0.07159928421353862; E=0.1249562713420228; F=16.16339869457496; G=0.001; H="False" }
    @{ Row=8; C=0.1059523809523822; D=0.0792738873881401; E=0.1326308745166243; F=17.42561354862084; G=0.001; H="False" }
    @{ Row=9; C=0.07110204081632582; D=0.04442354725208371; E=0.09778053438056793; F=11.69390130402446; G=0.001; H="False" }
    @{ Row=10; C=0.5364792989717934; D=0.5107061589429918; E=0.562252439000595; F=91.3322792732201; G=0.001; H="False" }
    @{ Row=11; C=0.5396696231932351; D=0.5139399515938827; E=0.5653992947925874; F=92.03062912083485; G=0.001; H="False" }
    @{ Row=12; C=0.5524502365386476; D=0.5267444061199772; E=0.5781560669573179; F=94.29750161967509; G=0.001; H="False" }
    @{ Row=13; C=0.5668387837062545; D=0.5411282024377084; E=0.5925493649748006; F=96.73559847418635; G=0.001; H="False" }
    @{ Row=14; C=0.4565128131573896; D=0.4308211835124244; E=0.4822044428023548; F=77.9650531778245; G=0.001; H="False" }
    @{ Row=15; C=0.4488382099827881; D=0.4231465803378229; E=0.4745298396277533; F=76.65435427216154; G=0.001; H="False" }
    @{ Row=16; C=0.4836885501188445; D=0.4579969204738794; E=0.5093801797638098; F=82.60623238743396; G=0.001; H="False" }
    @{ Row=17; C=0.003190324221441741; D=-0.01897452284434387; E=0.02535517128722736; F=0.631551458984317; G=0.9; H="True" }
    @{ Row=18; C=0.01597093756685422; D=-0.006166229423358961; E=0.0381081045570674; F=3.165534357324425; G=0.3826333963010845; H="True" }
    @{ Row=19; C=0.03035948473446113; D=0.008216801202295086; E=0.05250216826662717; F=6.015930412417551; G=0.001; H="False" }
    @{ Row=20; C=0.07996648581440371; D=0.05784581044828523; E=0.1020871611805222; F=15.86164733021868; G=0.001; H="False" }
    @{ Row=21; C=0.0876410889890052; D=0.06552041362288671; E=0.1097617643551237; F=17.38393316928175; G=0.001; H="False" }
    @{ Row=22; C=0.05279074885294881; D=0.03067007348683031; E=0.0749114242190673; F=10.47123969592766; G=0.001; H="False" }
    @{ Row=23; C=0.01278061334541248; D=-0.009305930431713207; E=0.03486715712253817; F=2.53899939006594; G=0.6612326889893446; H="True" }
    @{ Row=24; C=0.02716916051301939; D=0.005077087552974133; E=0.04926123347306464; F=5.396080394390007; G=0.004336209593780893; H="False" }
    @{ Row=25; C=0.08315681003584546; D=0.06108679571053278; E=0.1052268243611581; F=16.5323224302447; G=0.001; H="False" }
    @{ Row=26; C=0.09083141321044694; D=0.06876139888513426; E=0.1129014275357596; F=18.05810262974969; G=0.001; H="False" }
    @{ Row=27; C=0.05598107307439055; D=0.03391105874907788; E=0.07805108739970322; F=11.12954128059952; G=0.001; H="False" }
    @{ Row=28; C=0.01438854716760691; D=-0.007675754420589386; E=0.03645284875580319; F=2.861313345450331; G=0.522980241264126; H="True" }
    @{ Row=29; C=0.09593742338125794; D=0.07389520821983586; E=0.11797963854268; F=19.09727815295572; G=0.001; H="False" }
    @{ Row=30; C=0.1036120265558594; D=0.08156981139443734; E=0.1256542417172815; F=20.62498263336971; G=0.001; H="False" }
    @{ Row=31; C=0.06876168641980303; D=0.04671947125838095; E=0.09080390158122512; F=13.68768313285587; G=0.001; H="False" }
    @{ Row=32; C=0.1103259705488648; D=0.08827821508769604; E=0.1323737260100336; F=21.95593993637574; G=0.001; H="False" }
    @{ Row=33; C=0.1180005737234663; D=0.09595281826229753; E=0.1400483291846351; F=23.48326052552421; G=0.001; H="False" }
    @{ Row=34; C=0.08315023358740993; D=0.06110247812624113; E=0.1051979890485787; F=16.54770427360241; G=0.001; H="False" }
    @{ Row=35; C=0.007674603174601488; D=-0.01435104926801739; E=0.02970025561722036; F=1.528853274520744; G=0.9; H="True" }
    @{ Row=36; C=0.02717573696145491; D=0.005150084518836032; E=0.04920138940407379; F=5.413662895110165; G=0.004131024771766301; H="False" }
    @{ Row=37; C=0.0348503401360564; D=0.01282468769343752; E=0.05687599257867527; F=6.94251616963091; G=0.001; H="False" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G

    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = "'" + $item.H
    $hCell.Style = "Normal"
}

Write-Host "Updated rows 2-37 (C:H) with new Tukey HSD values"
